# ============================================================
# Edit script: theoretical_airdrop.xlsx - "v2 result in 0903"
#
# For each of the 5 sheets (TRX, JST, WBTT, WIN, NFT):
#   1. Append a new trailing date column AJ (header "20220211"),
#      mirroring the header style already used in row 1.
#   2. Insert a new period row "20210827-20210903" just above the
#      existing SUM row, holding the cohort constant across its
#      24-period diagonal run (columns M:AJ). The old SUM row
#      currently occupies row 13, so its leftover B:L values are
#      cleared before the new row is written there.
#   3. Push the old SUM row down one (row 13 -> row 14) and update
#      every column total to include the new row's contribution.
# ============================================================

$wb = $excel.ActiveWorkbook

# ---------------- TRX ----------------
$ws = $wb.Worksheets.Item("TRX")

# 1) New trailing column AJ: header + style copied from AI.
$ws.Range("AI1").Copy($ws.Range("AJ1"))
$ws.Range("AJ1").Value = "20220211"

# 2) Row 13 currently holds the old SUM row; grab its label style
#    for the new row-14 SUM label before overwriting row 13, then
#    clear the old SUM numbers out of B13:L13 (the new period row
#    only has data starting at column M).
$ws.Range("A13").Copy($ws.Range("A14"))
$ws.Range("A14").Value = "SUM"
$ws.Range("B13:AI13").ClearContents()

# New period row 13 ("20210827-20210903"); label style copied from A12.
$ws.Range("A12").Copy($ws.Range("A13"))
$ws.Range("A13").Value = "20210827-20210903"
$ws.Range("M13").Value = 837074.2362500001
$ws.Range("N13").Value = 837074.2362500001
$ws.Range("O13").Value = 837074.2362500001
$ws.Range("P13").Value = 837074.2362500001
$ws.Range("Q13").Value = 837074.2362500001
$ws.Range("R13").Value = 837074.2362500001
$ws.Range("S13").Value = 837074.2362500001
$ws.Range("T13").Value = 837074.2362500001
$ws.Range("U13").Value = 837074.2362500001
$ws.Range("V13").Value = 837074.2362500001
$ws.Range("W13").Value = 837074.2362500001
$ws.Range("X13").Value = 837074.2362500001
$ws.Range("Y13").Value = 837074.2362500001
$ws.Range("Z13").Value = 837074.2362500001
$ws.Range("AA13").Value = 837074.2362500001
$ws.Range("AB13").Value = 837074.2362500001
$ws.Range("AC13").Value = 837074.2362500001
$ws.Range("AD13").Value = 837074.2362500001
$ws.Range("AE13").Value = 837074.2362500001
$ws.Range("AF13").Value = 837074.2362500001
$ws.Range("AG13").Value = 837074.2362500001
$ws.Range("AH13").Value = 837074.2362500001
$ws.Range("AI13").Value = 837074.2362500001
$ws.Range("AJ13").Value = 837074.2362500001

# 3) New SUM row 14 column totals (old totals for B:L, + new row's
#    contribution for M:AJ).
$ws.Range("B14").Value = 939796.0891666667
$ws.Range("C14").Value = 1879592.178333333
$ws.Range("D14").Value = 2819388.2675
$ws.Range("E14").Value = 3759184.356666667
$ws.Range("F14").Value = 4801479.284583334
$ws.Range("G14").Value = 5843774.2125
$ws.Range("H14").Value = 6886069.140416667
$ws.Range("I14").Value = 7928364.068333333
$ws.Range("J14").Value = 8765438.304583333
$ws.Range("K14").Value = 9602512.540833334
$ws.Range("L14").Value = 10439586.77708333
$ws.Range("M14").Value = 11276661.01333333
$ws.Range("N14").Value = 11276661.01333333
$ws.Range("O14").Value = 11276661.01333333
$ws.Range("P14").Value = 11276661.01333333
$ws.Range("Q14").Value = 11276661.01333333
$ws.Range("R14").Value = 11276661.01333333
$ws.Range("S14").Value = 11276661.01333333
$ws.Range("T14").Value = 11276661.01333333
$ws.Range("U14").Value = 11276661.01333333
$ws.Range("V14").Value = 11276661.01333333
$ws.Range("W14").Value = 11276661.01333333
$ws.Range("X14").Value = 11276661.01333333
$ws.Range("Y14").Value = 11276661.01333333
$ws.Range("Z14").Value = 10336864.92416667
$ws.Range("AA14").Value = 9397068.834999999
$ws.Range("AB14").Value = 8457272.745833334
$ws.Range("AC14").Value = 7517476.656666666
$ws.Range("AD14").Value = 6475181.72875
$ws.Range("AE14").Value = 5432886.800833333
$ws.Range("AF14").Value = 4390591.872916667
$ws.Range("AG14").Value = 3348296.945
$ws.Range("AH14").Value = 2511222.70875
$ws.Range("AI14").Value = 1674148.4725
$ws.Range("AJ14").Value = 837074.2362500001

# ---------------- JST ----------------
$ws = $wb.Worksheets.Item("JST")

# 1) New trailing column AJ: header + style copied from AI.
$ws.Range("AI1").Copy($ws.Range("AJ1"))
$ws.Range("AJ1").Value = "20220211"

# 2) Row 13 currently holds the old SUM row; grab its label style
#    for the new row-14 SUM label before overwriting row 13, then
#    clear the old SUM numbers out of B13:L13 (the new period row
#    only has data starting at column M).
$ws.Range("A13").Copy($ws.Range("A14"))
$ws.Range("A14").Value = "SUM"
$ws.Range("B13:AI13").ClearContents()

# New period row 13 ("20210827-20210903"); label style copied from A12.
$ws.Range("A12").Copy($ws.Range("A13"))
$ws.Range("A13").Value = "20210827-20210903"
$ws.Range("M13").Value = 981750
$ws.Range("N13").Value = 981750
$ws.Range("O13").Value = 981750
$ws.Range("P13").Value = 981750
$ws.Range("Q13").Value = 981750
$ws.Range("R13").Value = 981750
$ws.Range("S13").Value = 981750
$ws.Range("T13").Value = 981750
$ws.Range("U13").Value = 981750
$ws.Range("V13").Value = 981750
$ws.Range("W13").Value = 981750
$ws.Range("X13").Value = 981750
$ws.Range("Y13").Value = 981750
$ws.Range("Z13").Value = 981750
$ws.Range("AA13").Value = 981750
$ws.Range("AB13").Value = 981750
$ws.Range("AC13").Value = 981750
$ws.Range("AD13").Value = 981750
$ws.Range("AE13").Value = 981750
$ws.Range("AF13").Value = 981750
$ws.Range("AG13").Value = 981750
$ws.Range("AH13").Value = 981750
$ws.Range("AI13").Value = 981750
$ws.Range("AJ13").Value = 981750

# 3) New SUM row 14 column totals (old totals for B:L, + new row's
#    contribution for M:AJ).
$ws.Range("B14").Value = 942083.3333333334
$ws.Range("C14").Value = 1829625
$ws.Range("D14").Value = 2717166.666666667
$ws.Range("E14").Value = 3604708.333333333
$ws.Range("F14").Value = 4477375
$ws.Range("G14").Value = 5350041.666666667
$ws.Range("H14").Value = 6222708.333333334
$ws.Range("I14").Value = 7095375.000000001
$ws.Range("J14").Value = 8022583.333333334
$ws.Range("K14").Value = 9004333.333333334
$ws.Range("L14").Value = 9986083.333333334
$ws.Range("M14").Value = 10967833.33333333
$ws.Range("N14").Value = 10967833.33333333
$ws.Range("O14").Value = 10967833.33333333
$ws.Range("P14").Value = 10967833.33333333
$ws.Range("Q14").Value = 10967833.33333333
$ws.Range("R14").Value = 10967833.33333333
$ws.Range("S14").Value = 10967833.33333333
$ws.Range("T14").Value = 10967833.33333333
$ws.Range("U14").Value = 10967833.33333333
$ws.Range("V14").Value = 10967833.33333333
$ws.Range("W14").Value = 10967833.33333333
$ws.Range("X14").Value = 10967833.33333333
$ws.Range("Y14").Value = 10967833.33333333
$ws.Range("Z14").Value = 10025750
$ws.Range("AA14").Value = 9138208.333333332
$ws.Range("AB14").Value = 8250666.666666666
$ws.Range("AC14").Value = 7363125
$ws.Range("AD14").Value = 6490458.333333334
$ws.Range("AE14").Value = 5617791.666666666
$ws.Range("AF14").Value = 4745125
$ws.Range("AG14").Value = 3872458.333333333
$ws.Range("AH14").Value = 2945250
$ws.Range("AI14").Value = 1963500
$ws.Range("AJ14").Value = 981750

# ---------------- WBTT ----------------
$ws = $wb.Worksheets.Item("WBTT")

# 1) New trailing column AJ: header + style copied from AI.
$ws.Range("AI1").Copy($ws.Range("AJ1"))
$ws.Range("AJ1").Value = "20220211"

# 2) Row 13 currently holds the old SUM row; grab its label style
#    for the new row-14 SUM label before overwriting row 13, then
#    clear the old SUM numbers out of B13:L13 (the new period row
#    only has data starting at column M).
$ws.Range("A13").Copy($ws.Range("A14"))
$ws.Range("A14").Value = "SUM"
$ws.Range("B13:AI13").ClearContents()

# New period row 13 ("20210827-20210903"); label style copied from A12.
$ws.Range("A12").Copy($ws.Range("A13"))
$ws.Range("A13").Value = "20210827-20210903"
$ws.Range("M13").Value = 2022976.932083334
$ws.Range("N13").Value = 2022976.932083334
$ws.Range("O13").Value = 2022976.932083334
$ws.Range("P13").Value = 2022976.932083334
$ws.Range("Q13").Value = 2022976.932083334
$ws.Range("R13").Value = 2022976.932083334
$ws.Range("S13").Value = 2022976.932083334
$ws.Range("T13").Value = 2022976.932083334
$ws.Range("U13").Value = 2022976.932083334
$ws.Range("V13").Value = 2022976.932083334
$ws.Range("W13").Value = 2022976.932083334
$ws.Range("X13").Value = 2022976.932083334
$ws.Range("Y13").Value = 2022976.932083334
$ws.Range("Z13").Value = 2022976.932083334
$ws.Range("AA13").Value = 2022976.932083334
$ws.Range("AB13").Value = 2022976.932083334
$ws.Range("AC13").Value = 2022976.932083334
$ws.Range("AD13").Value = 2022976.932083334
$ws.Range("AE13").Value = 2022976.932083334
$ws.Range("AF13").Value = 2022976.932083334
$ws.Range("AG13").Value = 2022976.932083334
$ws.Range("AH13").Value = 2022976.932083334
$ws.Range("AI13").Value = 2022976.932083334
$ws.Range("AJ13").Value = 2022976.932083334

# 3) New SUM row 14 column totals (old totals for B:L, + new row's
#    contribution for M:AJ).
$ws.Range("B14").Value = 2022976.932083334
$ws.Range("C14").Value = 4045953.864166667
$ws.Range("D14").Value = 6068930.796250001
$ws.Range("E14").Value = 8091907.728333334
$ws.Range("F14").Value = 10114884.66041667
$ws.Range("G14").Value = 12137861.5925
$ws.Range("H14").Value = 14160838.52458333
$ws.Range("I14").Value = 16183815.45666667
$ws.Range("J14").Value = 18206792.38875
$ws.Range("K14").Value = 20229769.32083334
$ws.Range("L14").Value = 22252746.25291667
$ws.Range("M14").Value = 24275723.18500001
$ws.Range("N14").Value = 24275723.18500001
$ws.Range("O14").Value = 24275723.18500001
$ws.Range("P14").Value = 24275723.18500001
$ws.Range("Q14").Value = 24275723.18500001
$ws.Range("R14").Value = 24275723.18500001
$ws.Range("S14").Value = 24275723.18500001
$ws.Range("T14").Value = 24275723.18500001
$ws.Range("U14").Value = 24275723.18500001
$ws.Range("V14").Value = 24275723.18500001
$ws.Range("W14").Value = 24275723.18500001
$ws.Range("X14").Value = 24275723.18500001
$ws.Range("Y14").Value = 24275723.18500001
$ws.Range("Z14").Value = 22252746.25291667
$ws.Range("AA14").Value = 20229769.32083334
$ws.Range("AB14").Value = 18206792.38875
$ws.Range("AC14").Value = 16183815.45666667
$ws.Range("AD14").Value = 14160838.52458333
$ws.Range("AE14").Value = 12137861.5925
$ws.Range("AF14").Value = 10114884.66041667
$ws.Range("AG14").Value = 8091907.728333334
$ws.Range("AH14").Value = 6068930.796250001
$ws.Range("AI14").Value = 4045953.864166667
$ws.Range("AJ14").Value = 2022976.932083334

# ---------------- WIN ----------------
$ws = $wb.Worksheets.Item("WIN")

# 1) New trailing column AJ: header + style copied from AI.
$ws.Range("AI1").Copy($ws.Range("AJ1"))
$ws.Range("AJ1").Value = "20220211"

# 2) Row 13 currently holds the old SUM row; grab its label style
#    for the new row-14 SUM label before overwriting row 13, then
#    clear the old SUM numbers out of B13:L13 (the new period row
#    only has data starting at column M).
$ws.Range("A13").Copy($ws.Range("A14"))
$ws.Range("A14").Value = "SUM"
$ws.Range("B13:AI13").ClearContents()

# New period row 13 ("20210827-20210903"); label style copied from A12.
$ws.Range("A12").Copy($ws.Range("A13"))
$ws.Range("A13").Value = "20210827-20210903"
$ws.Range("M13").Value = 5319315.755833333
$ws.Range("N13").Value = 5319315.755833333
$ws.Range("O13").Value = 5319315.755833333
$ws.Range("P13").Value = 5319315.755833333
$ws.Range("Q13").Value = 5319315.755833333
$ws.Range("R13").Value = 5319315.755833333
$ws.Range("S13").Value = 5319315.755833333
$ws.Range("T13").Value = 5319315.755833333
$ws.Range("U13").Value = 5319315.755833333
$ws.Range("V13").Value = 5319315.755833333
$ws.Range("W13").Value = 5319315.755833333
$ws.Range("X13").Value = 5319315.755833333
$ws.Range("Y13").Value = 5319315.755833333
$ws.Range("Z13").Value = 5319315.755833333
$ws.Range("AA13").Value = 5319315.755833333
$ws.Range("AB13").Value = 5319315.755833333
$ws.Range("AC13").Value = 5319315.755833333
$ws.Range("AD13").Value = 5319315.755833333
$ws.Range("AE13").Value = 5319315.755833333
$ws.Range("AF13").Value = 5319315.755833333
$ws.Range("AG13").Value = 5319315.755833333
$ws.Range("AH13").Value = 5319315.755833333
$ws.Range("AI13").Value = 5319315.755833333
$ws.Range("AJ13").Value = 5319315.755833333

# 3) New SUM row 14 column totals (old totals for B:L, + new row's
#    contribution for M:AJ).
$ws.Range("B14").Value = 5319315.755833333
$ws.Range("C14").Value = 10638631.51166667
$ws.Range("D14").Value = 15957947.2675
$ws.Range("E14").Value = 21277263.02333333
$ws.Range("F14").Value = 26596578.77916667
$ws.Range("G14").Value = 31915894.535
$ws.Range("H14").Value = 37235210.29083334
$ws.Range("I14").Value = 42554526.04666667
$ws.Range("J14").Value = 47873841.80250001
$ws.Range("K14").Value = 53193157.55833334
$ws.Range("L14").Value = 58512473.31416668
$ws.Range("M14").Value = 63831789.07000002
$ws.Range("N14").Value = 63831789.07000002
$ws.Range("O14").Value = 63831789.07000002
$ws.Range("P14").Value = 63831789.07000002
$ws.Range("Q14").Value = 63831789.07000002
$ws.Range("R14").Value = 63831789.07000002
$ws.Range("S14").Value = 63831789.07000002
$ws.Range("T14").Value = 63831789.07000002
$ws.Range("U14").Value = 63831789.07000002
$ws.Range("V14").Value = 63831789.07000002
$ws.Range("W14").Value = 63831789.07000002
$ws.Range("X14").Value = 63831789.07000002
$ws.Range("Y14").Value = 63831789.07000002
$ws.Range("Z14").Value = 58512473.31416668
$ws.Range("AA14").Value = 53193157.55833334
$ws.Range("AB14").Value = 47873841.80250001
$ws.Range("AC14").Value = 42554526.04666667
$ws.Range("AD14").Value = 37235210.29083334
$ws.Range("AE14").Value = 31915894.535
$ws.Range("AF14").Value = 26596578.77916667
$ws.Range("AG14").Value = 21277263.02333333
$ws.Range("AH14").Value = 15957947.2675
$ws.Range("AI14").Value = 10638631.51166667
$ws.Range("AJ14").Value = 5319315.755833333

# ---------------- NFT ----------------
$ws = $wb.Worksheets.Item("NFT")

# 1) New trailing column AJ: header + style copied from AI.
$ws.Range("AI1").Copy($ws.Range("AJ1"))
$ws.Range("AJ1").Value = "20220211"

# 2) Row 13 currently holds the old SUM row; grab its label style
#    for the new row-14 SUM label before overwriting row 13, then
#    clear the old SUM numbers out of B13:L13 (the new period row
#    only has data starting at column M).
$ws.Range("A13").Copy($ws.Range("A14"))
$ws.Range("A14").Value = "SUM"
$ws.Range("B13:AI13").ClearContents()

# New period row 13 ("20210827-20210903"); label style copied from A12.
$ws.Range("A12").Copy($ws.Range("A13"))
$ws.Range("A13").Value = "20210827-20210903"
$ws.Range("M13").Value = 2833161790.824167
$ws.Range("N13").Value = 2833161790.824167
$ws.Range("O13").Value = 2833161790.824167
$ws.Range("P13").Value = 2833161790.824167
$ws.Range("Q13").Value = 2833161790.824167
$ws.Range("R13").Value = 2833161790.824167
$ws.Range("S13").Value = 2833161790.824167
$ws.Range("T13").Value = 2833161790.824167
$ws.Range("U13").Value = 2833161790.824167
$ws.Range("V13").Value = 2833161790.824167
$ws.Range("W13").Value = 2833161790.824167
$ws.Range("X13").Value = 2833161790.824167
$ws.Range("Y13").Value = 2833161790.824167
$ws.Range("Z13").Value = 2833161790.824167
$ws.Range("AA13").Value = 2833161790.824167
$ws.Range("AB13").Value = 2833161790.824167
$ws.Range("AC13").Value = 2833161790.824167
$ws.Range("AD13").Value = 2833161790.824167
$ws.Range("AE13").Value = 2833161790.824167
$ws.Range("AF13").Value = 2833161790.824167
$ws.Range("AG13").Value = 2833161790.824167
$ws.Range("AH13").Value = 2833161790.824167
$ws.Range("AI13").Value = 2833161790.824167
$ws.Range("AJ13").Value = 2833161790.824167

# 3) New SUM row 14 column totals (old totals for B:L, + new row's
#    contribution for M:AJ).
$ws.Range("B14").Value = 2833161790.824167
$ws.Range("C14").Value = 5666323581.648334
$ws.Range("D14").Value = 8499485372.4725
$ws.Range("E14").Value = 11332647163.29667
$ws.Range("F14").Value = 14165808954.12083
$ws.Range("G14").Value = 16998970744.945
$ws.Range("H14").Value = 19832132535.76917
$ws.Range("I14").Value = 22665294326.59333
$ws.Range("J14").Value = 25498456117.4175
$ws.Range("K14").Value = 28331617908.24166
$ws.Range("L14").Value = 31164779699.06583
$ws.Range("M14").Value = 33997941489.89
$ws.Range("N14").Value = 33997941489.89
$ws.Range("O14").Value = 33997941489.89
$ws.Range("P14").Value = 33997941489.89
$ws.Range("Q14").Value = 33997941489.89
$ws.Range("R14").Value = 33997941489.89
$ws.Range("S14").Value = 33997941489.89
$ws.Range("T14").Value = 33997941489.89
$ws.Range("U14").Value = 33997941489.89
$ws.Range("V14").Value = 33997941489.89
$ws.Range("W14").Value = 33997941489.89
$ws.Range("X14").Value = 33997941489.89
$ws.Range("Y14").Value = 33997941489.89
$ws.Range("Z14").Value = 31164779699.06583
$ws.Range("AA14").Value = 28331617908.24166
$ws.Range("AB14").Value = 25498456117.4175
$ws.Range("AC14").Value = 22665294326.59333
$ws.Range("AD14").Value = 19832132535.76917
$ws.Range("AE14").Value = 16998970744.945
$ws.Range("AF14").Value = 14165808954.12083
$ws.Range("AG14").Value = 11332647163.29667
$ws.Range("AH14").Value = 8499485372.4725
$ws.Range("AI14").Value = 5666323581.648334
$ws.Range("AJ14").Value = 2833161790.824167

